$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 232, shifting existing rows 232:286 down to 233:287
$ws.Rows("232:232").Insert()

# Populate the newly inserted row with the new record
$ws.Range("A232").Value = 5
$ws.Range("B232").Value = "Macroferia Regional de Talca"
$ws.Range("C232").Value = "Maule"
$ws.Range("D232").Value = 44711
$ws.Range("E232").Value = 7
$ws.Range("F232").Value = 100112006
$ws.Range("G232").Value = "Repollo"
$ws.Range("H232").Value = "Crespo record"
$ws.Range("I232").Value = "Primera"
$ws.Range("J232").Value = 5000
$ws.Range("K232").Value = 1000
$ws.Range("L232").Value = 1000
$ws.Range("M232").Value = 1000
$ws.Range("N232").Value = "`$/unidad"
$ws.Range("O232").Value = "Región del Maule"
$ws.Range("P232").Value = 1000
$ws.Range("Q232").Value = 1
$ws.Range("R232").Value = "Hortaliza"
